# Auto-generated edit script: applies scheduled market-data refresh to Belias_Profits workbook
# Each sheet corresponds to a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For every touched leve row we rewrite the price/profit columns (H:N) with the refreshed values.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 174.6
$ws.Range("I5").Value = 138.44444
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 138.44444
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -23.44443999999999
$ws.Range("N5").Value = -730
# Row 11
$ws.Range("H11").Value = 36.6
$ws.Range("I11").Value = 36.6
$ws.Range("K11").Value = 36.6
$ws.Range("M11").Value = 103.4
# Row 70
$ws.Range("H70").Value = 1550.3
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1550.3
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4650.9
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5190.9
# Row 73
$ws.Range("H73").Value = 1550.3
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1550.3
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4650.9
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6522.9
# Row 80
$ws.Range("H80").Value = 2747.2
$ws.Range("I80").Value = 1971.7142
$ws.Range("J80").Value = 3425.75
$ws.Range("K80").Value = 5915.142599999999
$ws.Range("L80").Value = 10277.25
$ws.Range("M80").Value = -4917.142599999999
$ws.Range("N80").Value = -12273.25
# Row 82
$ws.Range("H82").Value = 2555.6924
$ws.Range("I82").Value = 1202.1818
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 3606.5454
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -3200.5454
$ws.Range("N82").Value = -30812
# Row 83
$ws.Range("H83").Value = 2747.2
$ws.Range("I83").Value = 1971.7142
$ws.Range("J83").Value = 3425.75
$ws.Range("K83").Value = 17745.4278
$ws.Range("L83").Value = 30831.75
$ws.Range("M83").Value = -12753.4278
$ws.Range("N83").Value = -40815.75
# Row 85
$ws.Range("H85").Value = 2555.6924
$ws.Range("I85").Value = 1202.1818
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 3606.5454
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -2202.5454
$ws.Range("N85").Value = -32808
# Row 88
$ws.Range("H88").Value = 11001.167
$ws.Range("I88").Value = 10001.5
$ws.Range("J88").Value = 11501
$ws.Range("K88").Value = 10001.5
$ws.Range("L88").Value = 11501
$ws.Range("M88").Value = -9595.5
$ws.Range("N88").Value = -12313
# Row 91
$ws.Range("H91").Value = 11001.167
$ws.Range("I91").Value = 10001.5
$ws.Range("J91").Value = 11501
$ws.Range("K91").Value = 10001.5
$ws.Range("L91").Value = 11501
$ws.Range("M91").Value = -8597.5
$ws.Range("N91").Value = -14309
# Row 126
$ws.Range("H126").Value = 22666.666
$ws.Range("J126").Value = 22666.666
$ws.Range("L126").Value = 22666.666
$ws.Range("N126").Value = -32546.666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 61
$ws.Range("H61").Value = 1330.7457
$ws.Range("I61").Value = 1180.5581
$ws.Range("J61").Value = 1734.375
$ws.Range("K61").Value = 1180.5581
$ws.Range("L61").Value = 1734.375
$ws.Range("M61").Value = -968.5581
$ws.Range("N61").Value = -2158.375
# Row 88
$ws.Range("H88").Value = 1896.9231
$ws.Range("I88").Value = 1308.5714
$ws.Range("J88").Value = 2583.3333
$ws.Range("K88").Value = 1308.5714
$ws.Range("L88").Value = 2583.3333
$ws.Range("M88").Value = -902.5714
$ws.Range("N88").Value = -3395.3333
# Row 91
$ws.Range("H91").Value = 1896.9231
$ws.Range("I91").Value = 1308.5714
$ws.Range("J91").Value = 2583.3333
$ws.Range("K91").Value = 1308.5714
$ws.Range("L91").Value = 2583.3333
$ws.Range("M91").Value = 95.42859999999996
$ws.Range("N91").Value = -5391.3333
# Row 136
$ws.Range("H136").Value = 1330.7457
$ws.Range("I136").Value = 1180.5581
$ws.Range("J136").Value = 1734.375
$ws.Range("K136").Value = 3541.6743
$ws.Range("L136").Value = 5203.125
$ws.Range("M136").Value = -991.6742999999997
$ws.Range("N136").Value = -10303.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1803.7916
$ws.Range("I86").Value = 1594.2632
$ws.Range("J86").Value = 2600
$ws.Range("K86").Value = 1594.2632
$ws.Range("L86").Value = 2600
$ws.Range("M86").Value = -471.2632000000001
$ws.Range("N86").Value = -4846
# Row 89
$ws.Range("H89").Value = 1803.7916
$ws.Range("I89").Value = 1594.2632
$ws.Range("J89").Value = 2600
$ws.Range("K89").Value = 7971.316000000001
$ws.Range("L89").Value = 13000
$ws.Range("M89").Value = -2355.316000000001
$ws.Range("N89").Value = -24232
# Row 134
$ws.Range("H134").Value = 788228.5600000001
$ws.Range("I134").Value = 1215820.2
$ws.Range("J134").Value = 4310.5
$ws.Range("K134").Value = 3647460.6
$ws.Range("L134").Value = 12931.5
$ws.Range("M134").Value = -3644925.6
$ws.Range("N134").Value = -18001.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2301.4285
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 2306.6667
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 2306.6667
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -3554.6667
# Row 65
$ws.Range("H65").Value = 2301.4285
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 2306.6667
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 11533.3335
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -17773.3335

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 1496574.1
$ws.Range("I121").Value = 650
$ws.Range("J121").Value = 1621234.5
$ws.Range("K121").Value = 1950
$ws.Range("L121").Value = 4863703.5
$ws.Range("M121").Value = -640
$ws.Range("N121").Value = -4866323.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2282.6155
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 2439.1428
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 2439.1428
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -4435.1428
# Row 83
$ws.Range("H83").Value = 2282.6155
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 2439.1428
$ws.Range("K83").Value = 10500
$ws.Range("L83").Value = 12195.714
$ws.Range("M83").Value = -5508
$ws.Range("N83").Value = -22179.714

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 624.44446
$ws.Range("I55").Value = 146.83333
$ws.Range("J55").Value = 863.25
$ws.Range("K55").Value = 146.83333
$ws.Range("L55").Value = 863.25
$ws.Range("M55").Value = 26.16667000000001
$ws.Range("N55").Value = -1209.25
# Row 68
$ws.Range("H68").Value = 12045.182
$ws.Range("I68").Value = 27075.5
$ws.Range("J68").Value = 3456.4285
$ws.Range("K68").Value = 27075.5
$ws.Range("L68").Value = 3456.4285
$ws.Range("M68").Value = -26326.5
$ws.Range("N68").Value = -4954.4285
# Row 71
$ws.Range("H71").Value = 12045.182
$ws.Range("I71").Value = 27075.5
$ws.Range("J71").Value = 3456.4285
$ws.Range("K71").Value = 135377.5
$ws.Range("L71").Value = 17282.1425
$ws.Range("M71").Value = -131633.5
$ws.Range("N71").Value = -24770.1425
# Row 82
$ws.Range("H82").Value = 1649.15
$ws.Range("I82").Value = 2925
$ws.Range("J82").Value = 1330.1875
$ws.Range("K82").Value = 2925
$ws.Range("L82").Value = 1330.1875
$ws.Range("M82").Value = -2564
$ws.Range("N82").Value = -2052.1875
# Row 85
$ws.Range("H85").Value = 1649.15
$ws.Range("I85").Value = 2925
$ws.Range("J85").Value = 1330.1875
$ws.Range("K85").Value = 2925
$ws.Range("L85").Value = 1330.1875
$ws.Range("M85").Value = -1677
$ws.Range("N85").Value = -3826.1875
# Row 136
$ws.Range("H136").Value = 1906.5405
$ws.Range("I136").Value = 1306
$ws.Range("J136").Value = 2787.3333
$ws.Range("K136").Value = 3918
$ws.Range("L136").Value = 8361.999899999999
$ws.Range("M136").Value = -1368
$ws.Range("N136").Value = -13461.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2371.2856
$ws.Range("I81").Value = 1079.8
$ws.Range("J81").Value = 5600
$ws.Range("K81").Value = 2159.6
$ws.Range("L81").Value = 11200
$ws.Range("M81").Value = -1098.6
$ws.Range("N81").Value = -13322
# Row 84
$ws.Range("H84").Value = 2371.2856
$ws.Range("I84").Value = 1079.8
$ws.Range("J84").Value = 5600
$ws.Range("K84").Value = 10798
$ws.Range("L84").Value = 56000
$ws.Range("M84").Value = -5494
$ws.Range("N84").Value = -66608
# Row 113
$ws.Range("H113").Value = 588.95654
$ws.Range("I113").Value = 440.94116
$ws.Range("J113").Value = 1008.3333
$ws.Range("K113").Value = 1322.82348
$ws.Range("L113").Value = 3024.9999
$ws.Range("M113").Value = 847.17652
$ws.Range("N113").Value = -7364.9999
# Row 136
$ws.Range("H136").Value = 1619.3334
$ws.Range("I136").Value = 1432.2881
$ws.Range("K136").Value = 4296.8643
$ws.Range("M136").Value = -1746.8643
